$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2, 3, 4, 8, 9, 10
$ws.Range("F2").Value = 14
$ws.Range("F3").Value = 4
$ws.Range("F4").Value = -3
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = -3
